$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_9_7_1"
$ws.Range("B2").Value = 0.008269583734393127
$ws.Range("C2").Value = -0.9980341123489331
$ws.Range("D2").Value = -1.825884951326852
$ws.Range("E2").Value = -0.8476703496255351
$ws.Range("F2").Value = 1.097553372383118
$ws.Range("G2").Value = 1.359334111213684
$ws.Range("H2").Value = 1.699893236160278
$ws.Range("I2").Value = 1.519587874412537
$ws.Range("A3").Value = "model_9_7_0"
$ws.Range("B3").Value = 0.01042201535442056
$ws.Range("C3").Value = -0.8426525177311099
$ws.Range("D3").Value = -1.804813901391807
$ws.Range("E3").Value = -0.772366664405618
$ws.Range("F3").Value = 1.095171213150024
$ws.Range("G3").Value = 1.253622531890869
$ws.Range("H3").Value = 1.687218070030212
$ws.Range("I3").Value = 1.457655549049377
$ws.Range("A4").Value = "model_9_7_2"
$ws.Range("B4").Value = 0.01775460137989338
$ws.Range("C4").Value = -1.099752366865607
$ws.Range("D4").Value = -1.855848913971139
$ws.Range("E4").Value = -0.9025508770059807
$ws.Range("F4").Value = 1.087056279182434
$ws.Range("G4").Value = 1.428536653518677
$ws.Range("H4").Value = 1.717917919158936
$ws.Range("I4").Value = 1.564723491668701
$ws.Range("A5").Value = "model_9_7_3"
$ws.Range("B5").Value = 0.0304512950525373
$ws.Range("C5").Value = -1.115351487856711
$ws.Range("D5").Value = -1.920042160240762
$ws.Range("E5").Value = -0.9314682215961558
$ws.Range("F5").Value = 1.073004841804504
$ws.Range("G5").Value = 1.439149379730225
$ws.Range("H5").Value = 1.756532907485962
$ws.Range("I5").Value = 1.588506102561951
$ws.Range("A6").Value = "model_9_7_4"
$ws.Range("B6").Value = 0.04471431979167506
$ws.Range("C6").Value = -1.168776620700449
$ws.Range("D6").Value = -1.983223901092988
$ws.Range("E6").Value = -0.9766129596924058
$ws.Range("F6").Value = 1.057219743728638
$ws.Range("G6").Value = 1.475496411323547
$ws.Range("H6").Value = 1.794539570808411
$ws.Range("I6").Value = 1.625634789466858
$ws.Range("A7").Value = "model_9_7_6"
$ws.Range("B7").Value = 0.05708626121222504
$ws.Range("C7").Value = -1.216055437063962
$ws.Range("D7").Value = -2.201672578133753
$ws.Range("E7").Value = -1.072506948749393
$ws.Range("F7").Value = 1.043527722358704
$ws.Range("G7").Value = 1.507661819458008
$ws.Range("H7").Value = 1.925945997238159
$ws.Range("I7").Value = 1.704501271247864
$ws.Range("B8").Value = 0.05719259986529035
$ws.Range("C8").Value = -1.163273067173712
$ws.Range("D8").Value = -2.034875755301694
$ws.Range("E8").Value = -0.9919725224352751
$ws.Range("F8").Value = 1.043410062789917
$ws.Range("G8").Value = 1.471752047538757
$ws.Range("H8").Value = 1.825610399246216
$ws.Range("I8").Value = 1.638267040252686
$ws.Range("A9").Value = "model_9_7_7"
$ws.Range("B9").Value = 0.07428783910005365
$ws.Range("C9").Value = -1.183169751859986
$ws.Range("D9").Value = -2.256437918703088
$ws.Range("E9").Value = -1.076948984393532
$ws.Range("F9").Value = 1.024490594863892
$ws.Range("G9").Value = 1.485288500785828
$ws.Range("H9").Value = 1.958889722824097
$ws.Range("I9").Value = 1.708154559135437
$ws.Range("A10").Value = "model_9_7_12"
$ws.Range("B10").Value = 0.1337159399209247
$ws.Range("C10").Value = -1.348043554010634
$ws.Range("D10").Value = -5.652616917718243
$ws.Range("E10").Value = -2.318114141950082
$ws.Range("F10").Value = 0.9587212204933167
$ws.Range("G10").Value = 1.597458004951477
$ws.Range("H10").Value = 4.001839637756348
$ws.Range("I10").Value = 2.728931665420532
$ws.Range("A11").Value = "model_9_7_13"
$ws.Range("B11").Value = 0.1368222699593887
$ws.Range("C11").Value = -1.314435477570808
$ws.Range("D11").Value = -5.67629933749576
$ws.Range("E11").Value = -2.311546431789816
$ws.Range("F11").Value = 0.9552834630012512
$ws.Range("G11").Value = 1.574593305587769
$ws.Range("H11").Value = 4.016086101531982
$ws.Range("I11").Value = 2.723530054092407
$ws.Range("A12").Value = "model_9_7_9"
$ws.Range("B12").Value = 0.1696650003697758
$ws.Range("C12").Value = -1.180301964355464
$ws.Range("D12").Value = -3.875359637561898
$ws.Range("E12").Value = -1.63292795400621
$ws.Range("F12").Value = 0.9189361929893494
$ws.Range("G12").Value = 1.48333740234375
$ws.Range("H12").Value = 2.932742118835449
$ws.Range("I12").Value = 2.165410995483398
$ws.Range("A13").Value = "model_9_7_11"
$ws.Range("B13").Value = 0.2015192028970284
$ws.Range("C13").Value = -1.302388846656332
$ws.Range("D13").Value = -4.359234302193378
$ws.Range("E13").Value = -1.852938779379683
$ws.Range("F13").Value = 0.883682906627655
$ws.Range("G13").Value = 1.566397547721863
$ws.Range("H13").Value = 3.223813533782959
$ws.Range("I13").Value = 2.346355438232422
$ws.Range("A14").Value = "model_9_7_10"
$ws.Range("B14").Value = 0.2088967129509747
$ws.Range("C14").Value = -1.314848225530678
$ws.Range("D14").Value = -4.142713617020108
$ws.Range("E14").Value = -1.78387846463489
$ws.Range("F14").Value = 0.8755182027816772
$ws.Range("G14").Value = 1.57487416267395
$ws.Range("H14").Value = 3.09356689453125
$ws.Range("I14").Value = 2.289557695388794
$ws.Range("A15").Value = "model_9_7_8"
$ws.Range("B15").Value = 0.2112965652143178
$ws.Range("C15").Value = -0.9121407461608129
$ws.Range("D15").Value = -2.286337507866506
$ws.Range("E15").Value = -0.9685405929453119
$ws.Range("F15").Value = 0.8728622794151306
$ws.Range("G15").Value = 1.300897836685181
$ws.Range("H15").Value = 1.97687554359436
$ws.Range("I15").Value = 1.618995785713196
$ws.Range("A16").Value = "model_9_7_22"
$ws.Range("B16").Value = 0.2480939639993316
$ws.Range("C16").Value = -1.655705585444055
$ws.Range("D16").Value = -5.078573859076574
$ws.Range("E16").Value = -2.25526894016945
$ws.Range("F16").Value = 0.8321383595466614
$ws.Range("G16").Value = 1.806771636009216
$ws.Range("H16").Value = 3.656527757644653
$ws.Range("I16").Value = 2.677245616912842
$ws.Range("A17").Value = "model_9_7_21"
$ws.Range("B17").Value = 0.252875290712947
$ws.Range("C17").Value = -1.568329496829891
$ws.Range("D17").Value = -5.101939113396917
$ws.Range("E17").Value = -2.22504086832386
$ws.Range("F17").Value = 0.8268469572067261
$ws.Range("G17").Value = 1.747326493263245
$ws.Range("H17").Value = 3.670583009719849
$ws.Range("I17").Value = 2.652384996414185
$ws.Range("A18").Value = "model_9_7_20"
$ws.Range("B18").Value = 0.2644003129939684
$ws.Range("C18").Value = -1.35320868641878
$ws.Range("D18").Value = -5.121259708978062
$ws.Range("E18").Value = -2.137489366317392
$ws.Range("F18").Value = 0.8140920996665955
$ws.Range("G18").Value = 1.600972056388855
$ws.Range("H18").Value = 3.682205200195312
$ws.Range("I18").Value = 2.580379724502563
$ws.Range("A19").Value = "model_9_7_19"
$ws.Range("B19").Value = 0.2708586017420108
$ws.Range("C19").Value = -1.353966679415374
$ws.Range("D19").Value = -5.027002185685125
$ws.Range("E19").Value = -2.105376678017563
$ws.Range("F19").Value = 0.8069446682929993
$ws.Range("G19").Value = 1.601487636566162
$ws.Range("H19").Value = 3.625504970550537
$ws.Range("I19").Value = 2.553969144821167
$ws.Range("A20").Value = "model_9_7_23"
$ws.Range("B20").Value = 0.2729026577035381
$ws.Range("C20").Value = -1.690860566002198
$ws.Range("D20").Value = -4.743082671955547
$ws.Range("E20").Value = -2.155195272862731
$ws.Range("F20").Value = 0.8046825528144836
$ws.Range("G20").Value = 1.830688953399658
$ws.Range("H20").Value = 3.454715251922607
$ws.Range("I20").Value = 2.59494161605835
$ws.Range("A21").Value = "model_9_7_24"
$ws.Range("B21").Value = 0.2743141064281475
$ws.Range("C21").Value = -1.683141423688768
$ws.Range("D21").Value = -4.742278719496897
$ws.Range("E21").Value = -2.151532205084019
$ws.Range("F21").Value = 0.8031204342842102
$ws.Range("G21").Value = 1.825437068939209
$ws.Range("H21").Value = 3.454231739044189
$ws.Range("I21").Value = 2.591928720474243
$ws.Range("A22").Value = "model_9_7_18"
$ws.Range("B22").Value = 0.3586639582477902
$ws.Range("C22").Value = -1.483140334658721
$ws.Range("D22").Value = -3.725530369354829
$ws.Range("E22").Value = -1.713983718903235
$ws.Range("F22").Value = 0.7097700238227844
$ws.Range("G22").Value = 1.689369201660156
$ws.Range("H22").Value = 2.842613220214844
$ws.Range("I22").Value = 2.232073783874512
$ws.Range("A23").Value = "model_9_7_14"
$ws.Range("B23").Value = 0.3891628890758279
$ws.Range("C23").Value = -1.499616368426848
$ws.Range("D23").Value = -2.951031342831101
$ws.Range("E23").Value = -1.454616696081315
$ws.Range("F23").Value = 0.6760166883468628
$ws.Range("G23").Value = 1.700578451156616
$ws.Range("H23").Value = 2.376718044281006
$ws.Range("I23").Value = 2.01876163482666
$ws.Range("A24").Value = "model_9_7_17"
$ws.Range("B24").Value = 0.4512284318815754
$ws.Range("C24").Value = -1.453145724825883
$ws.Range("D24").Value = -2.516106860447353
$ws.Range("E24").Value = -1.284565841456166
$ws.Range("F24").Value = 0.6073284149169922
$ws.Range("G24").Value = 1.668962955474854
$ws.Range("H24").Value = 2.115091800689697
$ws.Range("I24").Value = 1.878905773162842
$ws.Range("A25").Value = "model_9_7_16"
$ws.Range("B25").Value = 0.4775119151136697
$ws.Range("C25").Value = -1.467917005900941
$ws.Range("D25").Value = -2.148589897558486
$ws.Range("E25").Value = -1.164541047234125
$ws.Range("F25").Value = 0.5782403945922852
$ws.Range("G25").Value = 1.679012179374695
$ws.Range("H25").Value = 1.894014358520508
$ws.Range("I25").Value = 1.780193328857422
$ws.Range("A26").Value = "model_9_7_15"
$ws.Range("B26").Value = 0.478620731160034
$ws.Range("C26").Value = -1.478332272592475
$ws.Range("D26").Value = -2.090222764822247
$ws.Range("E26").Value = -1.149008292401647
$ws.Range("F26").Value = 0.5770131945610046
$ws.Range("G26").Value = 1.686098217964172
$ws.Range("H26").Value = 1.858904004096985
$ws.Range("I26").Value = 1.767418622970581
